# Commit: "Mecanismo de amacenamiento de imagenes" - Se implementa el mecanismo
# de almacenamiento de imagenes en el proyecto.
#
# Observable effect on "Plantilla Lista de Tareas de la 5ta Iteración.xlsx":
#   - On sheet "Casos de Uso", cell Q10 (Día 2 "Cons." for the image-storage
#     task row) is filled in with 1 hour consumed. This ripples through the
#     shared "restante" formulas for the rest of the row (R10..BA10).
#   - The view/selection state on that sheet moves on to cell Q11.
#   - The row-4 header mergeCells for the last few "day" columns get
#     re-created (re-merged), which re-orders them to the front of the
#     <mergeCells> list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- Data edit -----------------------------------------------------------
# Register 1 hour consumed on day 2 for the "Mecanismo para guardar, obtener
# y mostrar imagenes..." task row (row 10). The dependent shared formulas
# (R10:BA10) recompute automatically.
$ws.Range("Q10").Value = 1

# --- Merged header cells: re-create so they land at the front ------------
# (matches the reordering seen in the saved file: the last five day-header
# merges move ahead of the earlier ones in the <mergeCells> list)
$reorderedMerges = @("AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4")
$restOfMerges = @("AL4:AM4", "H4:I4", "K4:L4", "N4:O4", "Q4:R4", "T4:U4", "W4:X4", "Z4:AA4", "AC4:AD4", "AF4:AG4", "AI4:AJ4")
$allMerges = $reorderedMerges + $restOfMerges

foreach ($r in $allMerges) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $allMerges) {
    $ws.Range($r).Merge()
}

# --- View/selection state --------------------------------------------------
# Final active cell on the sheet ends up at Q11.
$ws.Range("Q11").Select()
